$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the submit logic columns
$ws.Range("G1").Value = "Standard"
$ws.Range("H1").Value = "Power Up"
$ws.Range("I1").Value = "Double"

# Update evaluate wins logic: loss row (row 12) now uses negative values
$ws.Range("G12").Value = -5
$ws.Range("H12").Value = -8

# Move the active selection to O7 to match the saved view state
$ws.Range("O7").Select()
